# Scheduled runner: refresh market-board snapshot values (currentAveragePrice*,
# LevePrice*, LeveProfit*) in the Shiva_Profits leve-crafting-profit sheets.
# Pricing pulled from the latest Universalis snapshot; downstream profit columns
# recompute from the refreshed prices.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
# ALC row 106: Enchanted Palladium Ink
$ws_ALC.Range("H106").Value = 125003860
$ws_ALC.Range("I106").Value = 125003860
$ws_ALC.Range("K106").Value = 125003860
$ws_ALC.Range("M106").Value = -125003229

# ALC row 132: Growth Formula Lambda
$ws_ALC.Range("H132").Value = 5442.885
$ws_ALC.Range("I132").Value = 4420.763
$ws_ALC.Range("J132").Value = 8217.214
$ws_ALC.Range("K132").Value = 13262.289
$ws_ALC.Range("L132").Value = 24651.642
$ws_ALC.Range("M132").Value = -10732.289
$ws_ALC.Range("N132").Value = -29711.642

# ALC row 135: Grade 1 Gemsap of Intelligence
$ws_ALC.Range("H135").Value = 953.931
$ws_ALC.Range("J135").Value = 3320.3333
$ws_ALC.Range("L135").Value = 29882.9997
$ws_ALC.Range("N135").Value = -34952.9997

# ALC row 137: Magnesia Whetstone
$ws_ALC.Range("H137").Value = 2490.1052
$ws_ALC.Range("I137").Value = 1880.5714
$ws_ALC.Range("J137").Value = 4196.8
$ws_ALC.Range("K137").Value = 5641.7142
$ws_ALC.Range("L137").Value = 12590.4
$ws_ALC.Range("M137").Value = -3091.7142
$ws_ALC.Range("N137").Value = -17690.4

$ws_ARM = $wb.Worksheets.Item("ARM")
# ARM row 5: Bronze Rivets
$ws_ARM.Range("H5").Value = 174.6
$ws_ARM.Range("I5").Value = 156
$ws_ARM.Range("K5").Value = 156
$ws_ARM.Range("M5").Value = -44

# ARM row 132: Mountain Chromite Ingot
$ws_ARM.Range("H132").Value = 2889.6365
$ws_ARM.Range("I132").Value = 2468.625
$ws_ARM.Range("J132").Value = 4012.3333
$ws_ARM.Range("K132").Value = 7405.875
$ws_ARM.Range("L132").Value = 12036.9999
$ws_ARM.Range("M132").Value = -4875.875
$ws_ARM.Range("N132").Value = -17096.9999

$ws_BSM = $wb.Worksheets.Item("BSM")
# BSM row 4: Bronze Rivets
$ws_BSM.Range("H4").Value = 174.6
$ws_BSM.Range("I4").Value = 156
$ws_BSM.Range("K4").Value = 156
$ws_BSM.Range("M4").Value = -41

# BSM row 20: Iron Ingot
$ws_BSM.Range("H20").Value = 1810.4865
$ws_BSM.Range("I20").Value = 1610.381
$ws_BSM.Range("J20").Value = 2073.125
$ws_BSM.Range("K20").Value = 1610.381
$ws_BSM.Range("L20").Value = 2073.125
$ws_BSM.Range("M20").Value = -1363.381
$ws_BSM.Range("N20").Value = -2567.125

# BSM row 107: Deepgold Nugget
$ws_BSM.Range("H107").Value = 7640.3335
$ws_BSM.Range("I107").Value = 5399.5
$ws_BSM.Range("J107").Value = 8088.5
$ws_BSM.Range("K107").Value = 5399.5
$ws_BSM.Range("L107").Value = 8088.5
$ws_BSM.Range("M107").Value = -3479.5
$ws_BSM.Range("N107").Value = -11928.5

# BSM row 122: High Durium Tathlums
$ws_BSM.Range("H122").Value = 200000
$ws_BSM.Range("J122").Value = 200000
$ws_BSM.Range("L122").Value = 200000
$ws_BSM.Range("N122").Value = -209800

# BSM row 134: Ruthenium Ingot
$ws_BSM.Range("H134").Value = 3130.577
$ws_BSM.Range("I134").Value = 2313.3171
$ws_BSM.Range("J134").Value = 6176.727
$ws_BSM.Range("K134").Value = 6939.951300000001
$ws_BSM.Range("L134").Value = 18530.181
$ws_BSM.Range("M134").Value = -4404.951300000001
$ws_BSM.Range("N134").Value = -23600.181

$ws_CRP = $wb.Worksheets.Item("CRP")
# CRP row 22: Elm Lumber
$ws_CRP.Range("H22").Value = 919.4286
$ws_CRP.Range("I22").Value = 732.1579
$ws_CRP.Range("J22").Value = 1141.8125
$ws_CRP.Range("K22").Value = 732.1579
$ws_CRP.Range("L22").Value = 1141.8125
$ws_CRP.Range("M22").Value = -382.1579
$ws_CRP.Range("N22").Value = -1841.8125

# CRP row 99: Pine Lumber
$ws_CRP.Range("H99").Value = 8794.637000000001
$ws_CRP.Range("I99").Value = 8148.278
$ws_CRP.Range("J99").Value = 9242.115
$ws_CRP.Range("K99").Value = 8148.278
$ws_CRP.Range("L99").Value = 9242.115
$ws_CRP.Range("M99").Value = -6650.278
$ws_CRP.Range("N99").Value = -12238.115

# CRP row 107: White Oak Lumber
$ws_CRP.Range("H107").Value = 1563.591
$ws_CRP.Range("I107").Value = 781.1429000000001
$ws_CRP.Range("J107").Value = 2932.875
$ws_CRP.Range("K107").Value = 781.1429000000001
$ws_CRP.Range("L107").Value = 2932.875
$ws_CRP.Range("M107").Value = 1138.8571
$ws_CRP.Range("N107").Value = -6772.875

# CRP row 126: Red Pine Lumber
$ws_CRP.Range("H126").Value = 8794.637000000001
$ws_CRP.Range("I126").Value = 8148.278
$ws_CRP.Range("J126").Value = 9242.115
$ws_CRP.Range("K126").Value = 24444.834
$ws_CRP.Range("L126").Value = 27726.345
$ws_CRP.Range("M126").Value = -21974.834
$ws_CRP.Range("N126").Value = -32666.345

# CRP row 132: Ginseng Lumber
$ws_CRP.Range("H132").Value = 8502.682000000001
$ws_CRP.Range("I132").Value = 9944.883
$ws_CRP.Range("J132").Value = 3599.2
$ws_CRP.Range("K132").Value = 29834.649
$ws_CRP.Range("L132").Value = 10797.6
$ws_CRP.Range("M132").Value = -27304.649
$ws_CRP.Range("N132").Value = -15857.6

# CRP row 134: Ceiba Lumber
$ws_CRP.Range("H134").Value = 6865.963
$ws_CRP.Range("I134").Value = 6005.5293
$ws_CRP.Range("J134").Value = 8328.700000000001
$ws_CRP.Range("K134").Value = 18016.5879
$ws_CRP.Range("L134").Value = 24986.1
$ws_CRP.Range("M134").Value = -15481.5879
$ws_CRP.Range("N134").Value = -30056.1

$ws_CUL = $wb.Worksheets.Item("CUL")
# CUL row 8: Sweet Cream
$ws_CUL.Range("H8").Value = 149.66667
$ws_CUL.Range("I8").Value = 149.66667
$ws_CUL.Range("K8").Value = 449.00001
$ws_CUL.Range("M8").Value = -310.00001

# CUL row 36: Crumpet
$ws_CUL.Range("H36").Value = 2294.4
$ws_CUL.Range("I36").Value = 367.25
$ws_CUL.Range("J36").Value = 10003
$ws_CUL.Range("K36").Value = 1101.75
$ws_CUL.Range("L36").Value = 30009
$ws_CUL.Range("M36").Value = -932.75
$ws_CUL.Range("N36").Value = -30347

# CUL row 99: Shorlog
$ws_CUL.Range("H99").Value = 2063.75
$ws_CUL.Range("I99").Value = 1305.3334
$ws_CUL.Range("J99").Value = 4339
$ws_CUL.Range("K99").Value = 3916.0002
$ws_CUL.Range("L99").Value = 13017
$ws_CUL.Range("M99").Value = -1670.0002
$ws_CUL.Range("N99").Value = -17509

$ws_GSM = $wb.Worksheets.Item("GSM")
# GSM row 102: Durium Ingot
$ws_GSM.Range("H102").Value = 6807.4287
$ws_GSM.Range("I102").Value = 10435.167
$ws_GSM.Range("K102").Value = 10435.167
$ws_GSM.Range("M102").Value = -8813.166999999999

# GSM row 107: Hard Mudstone Whetstone
$ws_GSM.Range("H107").Value = 1600.6666
$ws_GSM.Range("I107").Value = 2302
$ws_GSM.Range("J107").Value = 1250
$ws_GSM.Range("K107").Value = 2302
$ws_GSM.Range("L107").Value = 1250
$ws_GSM.Range("M107").Value = -382
$ws_GSM.Range("N107").Value = -5090

# GSM row 113: Manasilver Nugget
$ws_GSM.Range("H113").Value = 3262.25
$ws_GSM.Range("I113").Value = 2841.4
$ws_GSM.Range("K113").Value = 2841.4
$ws_GSM.Range("M113").Value = -671.4000000000001

# GSM row 139: White Gold Ring of Healing
$ws_GSM.Range("H139").Value = 150000
$ws_GSM.Range("J139").Value = 150000
$ws_GSM.Range("L139").Value = 150000
$ws_GSM.Range("N139").Value = -160280

$ws_LTW = $wb.Worksheets.Item("LTW")
# LTW row 7: Leather
$ws_LTW.Range("H7").Value = 3655.6775
$ws_LTW.Range("I7").Value = 3362.8262
$ws_LTW.Range("J7").Value = 4497.625
$ws_LTW.Range("K7").Value = 3362.8262
$ws_LTW.Range("L7").Value = 4497.625
$ws_LTW.Range("M7").Value = -3250.8262
$ws_LTW.Range("N7").Value = -4721.625

# LTW row 22: Aldgoat Leather
$ws_LTW.Range("H22").Value = 786.4211
$ws_LTW.Range("I22").Value = 783.5714
$ws_LTW.Range("J22").Value = 794.4
$ws_LTW.Range("K22").Value = 783.5714
$ws_LTW.Range("L22").Value = 794.4
$ws_LTW.Range("M22").Value = -488.5714
$ws_LTW.Range("N22").Value = -1384.4

# LTW row 27: Aldgoat Leather
$ws_LTW.Range("H27").Value = 786.4211
$ws_LTW.Range("I27").Value = 783.5714
$ws_LTW.Range("J27").Value = 794.4
$ws_LTW.Range("K27").Value = 783.5714
$ws_LTW.Range("L27").Value = 794.4
$ws_LTW.Range("M27").Value = -676.5714
$ws_LTW.Range("N27").Value = -1008.4

# LTW row 46: Boar Leather
$ws_LTW.Range("H46").Value = 2844.111
$ws_LTW.Range("I46").Value = 1149.8334
$ws_LTW.Range("J46").Value = 3328.1904
$ws_LTW.Range("K46").Value = 1149.8334
$ws_LTW.Range("L46").Value = 3328.1904
$ws_LTW.Range("M46").Value = -961.8334
$ws_LTW.Range("N46").Value = -3704.1904

# LTW row 126: Saiga Leather
$ws_LTW.Range("H126").Value = 3655.6775
$ws_LTW.Range("I126").Value = 3362.8262
$ws_LTW.Range("J126").Value = 4497.625
$ws_LTW.Range("K126").Value = 10088.4786
$ws_LTW.Range("L126").Value = 13492.875
$ws_LTW.Range("M126").Value = -7618.4786
$ws_LTW.Range("N126").Value = -18432.875

# LTW row 132: Silver Lobo Leather
$ws_LTW.Range("H132").Value = 59047.42
$ws_LTW.Range("I132").Value = 59047.42
$ws_LTW.Range("J132").Value = 0
$ws_LTW.Range("K132").Value = 177142.26
$ws_LTW.Range("L132").Value = 0
$ws_LTW.Range("M132").Value = -174612.26
$ws_LTW.Range("N132").ClearContents()

# LTW row 140: Gargantuaskin Shoes of Healing
$ws_LTW.Range("H140").Value = 88077.39999999999
$ws_LTW.Range("J140").Value = 88077.39999999999
$ws_LTW.Range("L140").Value = 88077.39999999999
$ws_LTW.Range("N140").Value = -98437.39999999999

$ws_WVR = $wb.Worksheets.Item("WVR")
# WVR row 42: Velveteen Gaskins
$ws_WVR.Range("H42").Value = 49999
$ws_WVR.Range("I42").Value = 49999
$ws_WVR.Range("K42").Value = 49999
$ws_WVR.Range("M42").Value = -49621

# WVR row 124: Almasty Serge Hat of Casting
$ws_WVR.Range("H124").Value = 166333
$ws_WVR.Range("J124").Value = 166333
$ws_WVR.Range("L124").Value = 166333
$ws_WVR.Range("N124").Value = -176153

# WVR row 126: Snow Linen
$ws_WVR.Range("H126").Value = 4799.6
$ws_WVR.Range("I126").Value = 4322.7666
$ws_WVR.Range("K126").Value = 12968.2998
$ws_WVR.Range("M126").Value = -10498.2998

# WVR row 132: Snow Cotton Cloth
$ws_WVR.Range("H132").Value = 2030.4546
$ws_WVR.Range("I132").Value = 1704.5555
$ws_WVR.Range("J132").Value = 3497
$ws_WVR.Range("K132").Value = 5113.666499999999
$ws_WVR.Range("L132").Value = 10491
$ws_WVR.Range("M132").Value = -2583.666499999999
$ws_WVR.Range("N132").Value = -15551
